# Auto-generated edit script applying the Ixion_Profits.xlsx value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 136.2
$ws.Range("I6").Value = 136.2
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 408.6
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -296.6
$ws.Range("N6").ClearContents()
# Row 8
$ws.Range("H8").Value = 57.375
$ws.Range("I8").Value = 57.375
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 172.125
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -33.125
$ws.Range("N8").ClearContents()
# Row 13
$ws.Range("H13").Value = 4266.6665
$ws.Range("J13").Value = 7000
$ws.Range("L13").Value = 7000
$ws.Range("N13").Value = -7338
# Row 137
$ws.Range("H137").Value = 1702.279
$ws.Range("I137").Value = 1194.8684
$ws.Range("J137").Value = 5558.6
$ws.Range("K137").Value = 3584.6052
$ws.Range("L137").Value = 16675.8
$ws.Range("M137").Value = -1034.6052
$ws.Range("N137").Value = -21775.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 3842.8572
$ws.Range("I3").Value = 1980
$ws.Range("K3").Value = 1980
$ws.Range("M3").Value = -1865
# Row 11
$ws.Range("H11").Value = 14500
$ws.Range("I11").Value = 14000
$ws.Range("K11").Value = 14000
$ws.Range("M11").Value = -13856
# Row 61
$ws.Range("H61").Value = 196711.44
$ws.Range("I61").Value = 5299.2334
$ws.Range("K61").Value = 5299.2334
$ws.Range("M61").Value = -5087.2334
# Row 74
$ws.Range("H74").Value = 1702.9459
$ws.Range("I74").Value = 1303.28
$ws.Range("J74").Value = 2535.5833
$ws.Range("K74").Value = 1303.28
$ws.Range("L74").Value = 2535.5833
$ws.Range("M74").Value = -429.28
$ws.Range("N74").Value = -4283.5833
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 77
$ws.Range("H77").Value = 1702.9459
$ws.Range("I77").Value = 1303.28
$ws.Range("J77").Value = 2535.5833
$ws.Range("K77").Value = 6516.4
$ws.Range("L77").Value = 12677.9165
$ws.Range("M77").Value = -2148.4
$ws.Range("N77").Value = -21413.9165
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 102
$ws.Range("H102").Value = 2470751.8
$ws.Range("I102").Value = 3368187.8
$ws.Range("J102").Value = 2803
$ws.Range("K102").Value = 3368187.8
$ws.Range("L102").Value = 2803
$ws.Range("M102").Value = -3366565.8
$ws.Range("N102").Value = -6047
# Row 132
$ws.Range("H132").Value = 2043510.9
$ws.Range("I132").Value = 1841.561
$ws.Range("J132").Value = 12507066
$ws.Range("K132").Value = 5524.683
$ws.Range("L132").Value = 37521198
$ws.Range("M132").Value = -2994.683
$ws.Range("N132").Value = -37526258
# Row 136
$ws.Range("H136").Value = 196711.44
$ws.Range("I136").Value = 5299.2334
$ws.Range("K136").Value = 15897.7002
$ws.Range("M136").Value = -13347.7002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 950.4167
$ws.Range("I11").Value = 540
$ws.Range("K11").Value = 540
$ws.Range("M11").Value = -400
# Row 105
$ws.Range("H105").Value = 2335.7144
$ws.Range("I105").Value = 2280
$ws.Range("K105").Value = 2280
$ws.Range("M105").Value = -533

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 200.1579
$ws.Range("I7").Value = 195.07143
$ws.Range("J7").Value = 214.4
$ws.Range("K7").Value = 195.07143
$ws.Range("L7").Value = 214.4
$ws.Range("M7").Value = -82.07142999999999
$ws.Range("N7").Value = -440.4
# Row 13
$ws.Range("H13").Value = 18840.715
$ws.Range("I13").Value = 185
$ws.Range("J13").Value = 26303
$ws.Range("K13").Value = 185
$ws.Range("L13").Value = 26303
$ws.Range("M13").Value = -46
$ws.Range("N13").Value = -26581
# Row 99
$ws.Range("H99").Value = 10503.833
$ws.Range("I99").Value = 10503.833
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 10503.833
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -9005.833000000001
$ws.Range("N99").ClearContents()
# Row 122
$ws.Range("H122").Value = 1352.1818
$ws.Range("I122").Value = 1418
$ws.Range("J122").Value = 694
$ws.Range("K122").Value = 4254
$ws.Range("L122").Value = 2082
$ws.Range("M122").Value = -1804
$ws.Range("N122").Value = -6982
# Row 126
$ws.Range("H126").Value = 10503.833
$ws.Range("I126").Value = 10503.833
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 31511.499
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -29041.499
$ws.Range("N126").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 900120
$ws.Range("I4").Value = 900120
$ws.Range("K4").Value = 2700360
$ws.Range("M4").Value = -2700248
# Row 11
$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 100
$ws.Range("K11").Value = 300
$ws.Range("M11").Value = -160
# Row 13
$ws.Range("H13").Value = 2339.375
$ws.Range("J13").Value = 2621.4285
$ws.Range("L13").Value = 7864.2855
$ws.Range("N13").Value = -8200.2855
# Row 137
$ws.Range("H137").Value = 17565.064
$ws.Range("I137").Value = 8155.5557
$ws.Range("J137").Value = 30593.615
$ws.Range("K137").Value = 24466.6671
$ws.Range("L137").Value = 91780.845
$ws.Range("M137").Value = -19366.6671
$ws.Range("N137").Value = -101980.845

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 12514000
$ws.Range("I7").Value = 12514000
$ws.Range("K7").Value = 12514000
$ws.Range("M7").Value = -12513888
# Row 8
$ws.Range("H8").Value = 12514000
$ws.Range("I8").Value = 12514000
$ws.Range("K8").Value = 12514000
$ws.Range("M8").Value = -12513861
# Row 9
$ws.Range("H9").Value = 90403.5
$ws.Range("I9").Value = 90403.5
$ws.Range("K9").Value = 90403.5
$ws.Range("M9").Value = -90233.5
# Row 11
$ws.Range("H11").Value = 15857143
$ws.Range("I11").Value = 13200000
$ws.Range("K11").Value = 13200000
$ws.Range("M11").Value = -13199861
# Row 70
$ws.Range("H70").Value = 5170.75
$ws.Range("I70").Value = 5399.923
$ws.Range("J70").Value = 4899.909
$ws.Range("K70").Value = 5399.923
$ws.Range("L70").Value = 4899.909
$ws.Range("M70").Value = -5129.923
$ws.Range("N70").Value = -5439.909
# Row 73
$ws.Range("H73").Value = 5170.75
$ws.Range("I73").Value = 5399.923
$ws.Range("J73").Value = 4899.909
$ws.Range("K73").Value = 5399.923
$ws.Range("L73").Value = 4899.909
$ws.Range("M73").Value = -4463.923
$ws.Range("N73").Value = -6771.909
# Row 107
$ws.Range("H107").Value = 1816
$ws.Range("I107").Value = 655.5
$ws.Range("J107").Value = 4601.2
$ws.Range("K107").Value = 655.5
$ws.Range("L107").Value = 4601.2
$ws.Range("M107").Value = 1264.5
$ws.Range("N107").Value = -8441.200000000001
# Row 113
$ws.Range("H113").Value = 34484036
$ws.Range("I113").Value = 40000960
$ws.Range("J113").Value = 3275
$ws.Range("K113").Value = 40000960
$ws.Range("L113").Value = 3275
$ws.Range("M113").Value = -39998790
$ws.Range("N113").Value = -7615
# Row 132
$ws.Range("H132").Value = 4876.559
$ws.Range("I132").Value = 7476.0386
$ws.Range("J132").Value = 2828.4849
$ws.Range("K132").Value = 22428.1158
$ws.Range("L132").Value = 8485.4547
$ws.Range("M132").Value = -19898.1158
$ws.Range("N132").Value = -13545.4547

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 63
$ws.Range("H63").Value = 32300
$ws.Range("J63").Value = 32300
$ws.Range("L63").Value = 32300
$ws.Range("N63").Value = -33798
# Row 66
$ws.Range("H66").Value = 32300
$ws.Range("J66").Value = 32300
$ws.Range("L66").Value = 96900
$ws.Range("N66").Value = -104388
# Row 68
$ws.Range("H68").Value = 1701.8948
$ws.Range("I68").Value = 1498.6666
$ws.Range("J68").Value = 2464
$ws.Range("K68").Value = 1498.6666
$ws.Range("L68").Value = 2464
$ws.Range("M68").Value = -749.6666
$ws.Range("N68").Value = -3962
# Row 71
$ws.Range("H71").Value = 1701.8948
$ws.Range("I71").Value = 1498.6666
$ws.Range("J71").Value = 2464
$ws.Range("K71").Value = 7493.333000000001
$ws.Range("L71").Value = 12320
$ws.Range("M71").Value = -3749.333000000001
$ws.Range("N71").Value = -19808
# Row 100
$ws.Range("H100").Value = 1626.0834
$ws.Range("I100").Value = 1531.625
$ws.Range("J100").Value = 1815
$ws.Range("K100").Value = 1531.625
$ws.Range("L100").Value = 1815
$ws.Range("M100").Value = -990.625
$ws.Range("N100").Value = -2897
# Row 122
$ws.Range("H122").Value = 8151937
$ws.Range("I122").Value = 10215910
$ws.Range("J122").Value = 3336000
$ws.Range("K122").Value = 30647730
$ws.Range("L122").Value = 10008000
$ws.Range("M122").Value = -30645280
$ws.Range("N122").Value = -10012900

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 1000
$ws.Range("J8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("N8").Value = -1280
# Row 136
$ws.Range("H136").Value = 2385.5344
$ws.Range("I136").Value = 2703.6553
$ws.Range("J136").Value = 2067.4138
$ws.Range("K136").Value = 8110.965899999999
$ws.Range("L136").Value = 6202.241399999999
$ws.Range("M136").Value = -5560.965899999999
$ws.Range("N136").Value = -11302.2414
